# Generate Report for Handoff
# Updates the localization-status report: marks the handoff type ("ht")
# for the "Ready for handoff" rows and refreshes the "Latest Handoff
# Datetime" / "Latest HO Xliff Generate Date" timestamps produced by the
# report run.

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 10, 12, 13, 14)

# zh-cn sheet: set Priority ("ht") and refresh the handoff datetime.
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZh.Range("E$r").Value = "ht"
    $wsZh.Range("H$r").Value = "2016-08-14 00:28:55"
}

# de-de sheet: set Priority ("ht") and refresh the handoff datetime.
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDe.Range("E$r").Value = "ht"
    $wsDe.Range("H$r").Value = "2016-08-14 00:29:06"
}

# Overview sheet: refresh the "Latest HO Xliff Generate Date" column.
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-14 00:29:06"
}
